$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2047
$ws.Range("J40").Value = 2085.4546
$ws.Range("L40").Value = 2085.4546
$ws.Range("N40").Value = -2435.4546

$ws.Range("H62").Value = 4871.4287
$ws.Range("I62").Value = 6700
$ws.Range("J62").Value = 2433.3333
$ws.Range("K62").Value = 6700
$ws.Range("L62").Value = 2433.3333
$ws.Range("M62").Value = -6076
$ws.Range("N62").Value = -3681.3333

$ws.Range("H64").Value = 2484590.5
$ws.Range("I64").Value = 5497488.5
$ws.Range("J64").Value = 3380.4707
$ws.Range("K64").Value = 5497488.5
$ws.Range("L64").Value = 3380.4707
$ws.Range("M64").Value = -5497240.5
$ws.Range("N64").Value = -3876.4707

$ws.Range("H65").Value = 4871.4287
$ws.Range("I65").Value = 6700
$ws.Range("J65").Value = 2433.3333
$ws.Range("K65").Value = 33500
$ws.Range("L65").Value = 12166.6665
$ws.Range("M65").Value = -30380
$ws.Range("N65").Value = -18406.6665

$ws.Range("H67").Value = 2484590.5
$ws.Range("I67").Value = 5497488.5
$ws.Range("J67").Value = 3380.4707
$ws.Range("K67").Value = 5497488.5
$ws.Range("L67").Value = 3380.4707
$ws.Range("M67").Value = -5496630.5
$ws.Range("N67").Value = -5096.4707

$ws.Range("H100").Value = 2415.125
$ws.Range("I100").Value = 1902.5
$ws.Range("J100").Value = 2586
$ws.Range("K100").Value = 1902.5
$ws.Range("L100").Value = 2586
$ws.Range("M100").Value = -1361.5
$ws.Range("N100").Value = -3668

$ws.Range("H106").Value = 716915.7
$ws.Range("I106").Value = 1001882
$ws.Range("K106").Value = 1001882
$ws.Range("M106").Value = -1001251

$ws.Range("H116").Value = 3250
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -10384

$ws.Range("H129").Value = 1574.2759
$ws.Range("I129").Value = 367.33334
$ws.Range("J129").Value = 1889.1305
$ws.Range("K129").Value = 1102.00002
$ws.Range("L129").Value = 5667.3915
$ws.Range("M129").Value = 3897.99998
$ws.Range("N129").Value = -15667.3915

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1078.8864
$ws.Range("I74").Value = 1108.4193
$ws.Range("J74").Value = 1008.46155
$ws.Range("K74").Value = 1108.4193
$ws.Range("L74").Value = 1008.46155
$ws.Range("M74").Value = -234.4193
$ws.Range("N74").Value = -2756.46155

$ws.Range("H77").Value = 1078.8864
$ws.Range("I77").Value = 1108.4193
$ws.Range("J77").Value = 1008.46155
$ws.Range("K77").Value = 5542.0965
$ws.Range("L77").Value = 5042.30775
$ws.Range("M77").Value = -1174.0965
$ws.Range("N77").Value = -13778.30775

$ws.Range("H97").Value = 25880.25
$ws.Range("I97").Value = 33503.332
$ws.Range("J97").Value = 3011
$ws.Range("K97").Value = 33503.332
$ws.Range("L97").Value = 3011
$ws.Range("M97").Value = -33007.332
$ws.Range("N97").Value = -4003

$ws.Range("H132").Value = 31253968
$ws.Range("I132").Value = 125005010
$ws.Range("J132").Value = 3623.3333
$ws.Range("K132").Value = 375015030
$ws.Range("L132").Value = 10869.9999
$ws.Range("M132").Value = -375012500
$ws.Range("N132").Value = -15929.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3361.342
$ws.Range("J105").Value = 4595.4546
$ws.Range("L105").Value = 4595.4546
$ws.Range("N105").Value = -8089.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3941.4
$ws.Range("I86").Value = 4178.5
$ws.Range("J86").Value = 3783.3333
$ws.Range("K86").Value = 4178.5
$ws.Range("L86").Value = 3783.3333
$ws.Range("M86").Value = -3055.5
$ws.Range("N86").Value = -6029.3333

$ws.Range("H89").Value = 3941.4
$ws.Range("I89").Value = 4178.5
$ws.Range("J89").Value = 3783.3333
$ws.Range("K89").Value = 20892.5
$ws.Range("L89").Value = 18916.6665
$ws.Range("M89").Value = -15276.5
$ws.Range("N89").Value = -30148.6665

$ws.Range("H99").Value = 2043.8889
$ws.Range("I99").Value = 1999.375
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 1999.375
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -501.375
$ws.Range("N99").Value = -5396

$ws.Range("H126").Value = 2043.8889
$ws.Range("I126").Value = 1999.375
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 5998.125
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -3528.125
$ws.Range("N126").Value = -12140

$ws.Range("H132").Value = 3884.353
$ws.Range("I132").Value = 3772
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 11316
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -8786
$ws.Range("N132").Value = -17808.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1984292.6
$ws.Range("I2").Value = 193.33333
$ws.Range("J2").Value = 13888888
$ws.Range("K2").Value = 1159.99998
$ws.Range("L2").Value = 83333328
$ws.Range("M2").Value = -1046.99998
$ws.Range("N2").Value = -83333554

$ws.Range("H21").Value = 283.33334
$ws.Range("I21").Value = 250
$ws.Range("J21").Value = 1250
$ws.Range("K21").Value = 750
$ws.Range("L21").Value = 3750
$ws.Range("M21").Value = -577
$ws.Range("N21").Value = -4096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4508.3687
$ws.Range("I126").Value = 1502.4
$ws.Range("J126").Value = 5581.9287
$ws.Range("K126").Value = 4507.200000000001
$ws.Range("L126").Value = 16745.7861
$ws.Range("M126").Value = -2037.200000000001
$ws.Range("N126").Value = -21685.7861

$ws.Range("H132").Value = 5318.1333
$ws.Range("I132").Value = 5597.1113
$ws.Range("J132").Value = 4899.6665
$ws.Range("K132").Value = 16791.3339
$ws.Range("L132").Value = 14698.9995
$ws.Range("M132").Value = -14261.3339
$ws.Range("N132").Value = -19758.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6071.8066
$ws.Range("I122").Value = 6385.684
$ws.Range("K122").Value = 19157.052
$ws.Range("M122").Value = -16707.052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1236.909
$ws.Range("I132").Value = 1206.2258
$ws.Range("J132").Value = 1310.0769
$ws.Range("K132").Value = 3618.6774
$ws.Range("L132").Value = 3930.2307
$ws.Range("M132").Value = -1088.6774
$ws.Range("N132").Value = -8990.2307
